$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 21.1307724478878
$ws.Cells.Item(2, 3).Value = 10.93865540096599
$ws.Cells.Item(2, 5).Value = 10.1068152764238
$ws.Cells.Item(2, 6).Value = 49.6830764657261
$ws.Cells.Item(2, 7).Value = 3.750388757801944
$ws.Cells.Item(2, 9).Value = 36.57626929566658
$ws.Cells.Item(2, 10).Value = 10.2912964837885
$ws.Cells.Item(2, 12).Value = 11.39285899582205
$ws.Cells.Item(2, 13).Value = 18.57488265708863

$ws.Cells.Item(3, 2).Value = 20.88313522863472
$ws.Cells.Item(3, 3).Value = 10.44892478237944
$ws.Cells.Item(3, 5).Value = 10.05357985249676
$ws.Cells.Item(3, 6).Value = 49.50226214171023
$ws.Cells.Item(3, 7).Value = 3.754635444520498
$ws.Cells.Item(3, 9).Value = 36.4820632487174
$ws.Cells.Item(3, 10).Value = 10.30676921505165
$ws.Cells.Item(3, 12).Value = 11.41594134071745
$ws.Cells.Item(3, 13).Value = 18.5641752541184

$ws.Cells.Item(4, 2).Value = 20.7365304265991
$ws.Cells.Item(4, 3).Value = 10.13996049795043
$ws.Cells.Item(4, 5).Value = 10.02019853670067
$ws.Cells.Item(4, 6).Value = 49.40107105256259
$ws.Cells.Item(4, 7).Value = 3.757376674370179
$ws.Cells.Item(4, 9).Value = 36.43061912813563
$ws.Cells.Item(4, 10).Value = 10.31682599618265
$ws.Cells.Item(4, 12).Value = 11.43156158228478
$ws.Cells.Item(4, 13).Value = 18.56171826021271

$ws.Cells.Item(5, 2).Value = 20.67823125994583
$ws.Cells.Item(5, 3).Value = 10.01218485203055
$ws.Cells.Item(5, 5).Value = 10.00642132125682
$ws.Cells.Item(5, 6).Value = 49.36232404868013
$ws.Cells.Item(5, 7).Value = 3.758527513439943
$ws.Cells.Item(5, 9).Value = 36.41126791878658
$ws.Cells.Item(5, 10).Value = 10.32106414938
$ws.Cells.Item(5, 12).Value = 11.43829168267511
$ws.Cells.Item(5, 13).Value = 18.56175495923125

$ws.Cells.Item(6, 2).Value = 20.66864001519183
$ws.Cells.Item(6, 3).Value = 9.990861473013609
$ws.Cells.Item(6, 5).Value = 10.00412312918891
$ws.Cells.Item(6, 6).Value = 49.35604092997938
$ws.Cells.Item(6, 7).Value = 3.758720652681392
$ws.Cells.Item(6, 9).Value = 36.40815214458974
$ws.Cells.Item(6, 10).Value = 10.32177634334921
$ws.Cells.Item(6, 12).Value = 11.43943126007769
$ws.Cells.Item(6, 13).Value = 18.56182380599275

$ws.Cells.Item(7, 2).Value = 20.73573824351788
$ws.Cells.Item(7, 3).Value = 10.13824454769821
$ws.Cells.Item(7, 5).Value = 10.02001343628026
$ws.Cells.Item(7, 6).Value = 49.40053839771874
$ws.Cells.Item(7, 7).Value = 3.757392058097003
$ws.Cells.Item(7, 9).Value = 36.43035161774809
$ws.Cells.Item(7, 10).Value = 10.3168825868314
$ws.Cells.Item(7, 12).Value = 11.43165086906541
$ws.Cells.Item(7, 13).Value = 18.56171454964582

$ws.Cells.Item(8, 2).Value = 21.0443023704971
$ws.Cells.Item(8, 3).Value = 10.77162940697876
$ws.Cells.Item(8, 5).Value = 10.08860421012911
$ws.Cells.Item(8, 6).Value = 49.6187028839329
$ws.Cells.Item(8, 7).Value = 3.751825336956548
$ws.Cells.Item(8, 9).Value = 36.54246161844855
$ws.Cells.Item(8, 10).Value = 10.29651599236082
$ws.Cells.Item(8, 12).Value = 11.40051777414732
$ws.Cells.Item(8, 13).Value = 18.57033791193368

$ws.Cells.Item(9, 2).Value = 21.68920412804532
$ws.Cells.Item(9, 3).Value = 11.94022035571664
$ws.Cells.Item(9, 5).Value = 10.21757646119735
$ws.Cells.Item(9, 6).Value = 50.12362385239805
$ws.Cells.Item(9, 7).Value = 3.741964127370006
$ws.Cells.Item(9, 9).Value = 36.81286004613271
$ws.Cells.Item(9, 10).Value = 10.26099180581486
$ws.Cells.Item(9, 12).Value = 11.35092021458638
$ws.Cells.Item(9, 13).Value = 18.61977958822673

$ws.Cells.Item(10, 2).Value = 22.18258559501976
$ws.Cells.Item(10, 3).Value = 12.74510620935807
$ws.Cells.Item(10, 5).Value = 10.30893236789282
$ws.Cells.Item(10, 6).Value = 50.54021924767995
$ws.Cells.Item(10, 7).Value = 3.735353767964646
$ws.Cells.Item(10, 9).Value = 37.04195898466156
$ws.Cells.Item(10, 10).Value = 10.23758213610583
$ws.Cells.Item(10, 12).Value = 11.32142106817061
$ws.Cells.Item(10, 13).Value = 18.67571492515855

$ws.Cells.Item(11, 2).Value = 22.41026764943875
$ws.Cells.Item(11, 3).Value = 13.09815844135032
$ws.Cells.Item(11, 5).Value = 10.34974855186563
$ws.Cells.Item(11, 6).Value = 50.73929777125068
$ws.Cells.Item(11, 7).Value = 3.732482526652642
$ws.Cells.Item(11, 9).Value = 37.15267949151917
$ws.Cells.Item(11, 10).Value = 10.22751635545467
$ws.Cells.Item(11, 12).Value = 11.30949890241647
$ws.Cells.Item(11, 13).Value = 18.70535950959064

$ws.Cells.Item(12, 2).Value = 22.49686553733859
$ws.Cells.Item(12, 3).Value = 13.22986987443831
$ws.Cells.Item(12, 5).Value = 10.36509719441714
$ws.Cells.Item(12, 6).Value = 50.81602552219493
$ws.Cells.Item(12, 7).Value = 3.731414655381653
$ws.Cells.Item(12, 9).Value = 37.19552939565334
$ws.Cells.Item(12, 10).Value = 10.22378861624723
$ws.Cells.Item(12, 12).Value = 11.30519880065348
$ws.Cells.Item(12, 13).Value = 18.71718284524942

$ws.Cells.Item(13, 2).Value = 22.47819961307081
$ws.Cells.Item(13, 3).Value = 13.20159301043385
$ws.Cells.Item(13, 5).Value = 10.36179638708589
$ws.Cells.Item(13, 6).Value = 50.7994417635439
$ws.Cells.Item(13, 7).Value = 3.731643779547372
$ws.Cells.Item(13, 9).Value = 37.18626009096148
$ws.Cells.Item(13, 10).Value = 10.22458771668103
$ws.Cells.Item(13, 12).Value = 11.30611537338471
$ws.Cells.Item(13, 13).Value = 18.71461000766203

$ws.Cells.Item(14, 2).Value = 22.41738499260534
$ws.Cells.Item(14, 3).Value = 13.10903464904117
$ws.Cells.Item(14, 5).Value = 10.35101345450552
$ws.Cells.Item(14, 6).Value = 50.74558354057548
$ws.Cells.Item(14, 7).Value = 3.732394284062091
$ws.Cells.Item(14, 9).Value = 37.15618636865647
$ws.Cells.Item(14, 10).Value = 10.2272079891823
$ws.Cells.Item(14, 12).Value = 11.3091408344689
$ws.Cells.Item(14, 13).Value = 18.7063202844274

$ws.Cells.Item(15, 2).Value = 22.38018117419667
$ws.Cells.Item(15, 3).Value = 13.05207922636501
$ws.Cells.Item(15, 5).Value = 10.34439456396952
$ws.Cells.Item(15, 6).Value = 50.71276739465054
$ws.Cells.Item(15, 7).Value = 3.732856513100879
$ws.Cells.Item(15, 9).Value = 37.13788507852627
$ws.Cells.Item(15, 10).Value = 10.22882391689093
$ws.Cells.Item(15, 12).Value = 11.31102193758071
$ws.Cells.Item(15, 13).Value = 18.70132021444294

$ws.Cells.Item(16, 2).Value = 22.16776384432449
$ws.Cells.Item(16, 3).Value = 12.7217611134928
$ws.Cells.Item(16, 5).Value = 10.30624992935023
$ws.Cells.Item(16, 6).Value = 50.52739860237769
$ws.Cells.Item(16, 7).Value = 3.735544132857633
$ws.Cells.Item(16, 9).Value = 37.0348529396888
$ws.Cells.Item(16, 10).Value = 10.23825170312388
$ws.Cells.Item(16, 12).Value = 11.32223027450167
$ws.Cells.Item(16, 13).Value = 18.67386156008947

$ws.Cells.Item(17, 2).Value = 22.03822097777795
$ws.Cells.Item(17, 3).Value = 12.51569245437779
$ws.Cells.Item(17, 5).Value = 10.28265876659124
$ws.Cells.Item(17, 6).Value = 50.4161073195005
$ws.Cells.Item(17, 7).Value = 3.737227601793621
$ws.Cells.Item(17, 9).Value = 36.97330295897346
$ws.Cells.Item(17, 10).Value = 10.24418480856636
$ws.Cells.Item(17, 12).Value = 11.32948916936025
$ws.Cells.Item(17, 13).Value = 18.65808783985169

$ws.Cells.Item(18, 2).Value = 21.96402148138907
$ws.Cells.Item(18, 3).Value = 12.39593953590017
$ws.Cells.Item(18, 5).Value = 10.2690199577757
$ws.Cells.Item(18, 6).Value = 50.35299826474018
$ws.Cells.Item(18, 7).Value = 3.738208682589045
$ws.Cells.Item(18, 9).Value = 36.93851390556014
$ws.Cells.Item(18, 10).Value = 10.24765228183874
$ws.Cells.Item(18, 12).Value = 11.33380525601731
$ws.Cells.Item(18, 13).Value = 18.6494108889842

$ws.Cells.Item(19, 2).Value = 21.93895471513758
$ws.Cells.Item(19, 3).Value = 12.35518567330468
$ws.Cells.Item(19, 5).Value = 10.2643901296983
$ws.Cells.Item(19, 6).Value = 50.33178668840787
$ws.Cells.Item(19, 7).Value = 3.738543060995905
$ws.Cells.Item(19, 9).Value = 36.92684054962983
$ws.Cells.Item(19, 10).Value = 10.2488357373853
$ws.Cells.Item(19, 12).Value = 11.33529084106345
$ws.Cells.Item(19, 13).Value = 18.64654116135699

$ws.Cells.Item(20, 2).Value = 22.05197956726895
$ws.Cells.Item(20, 3).Value = 12.53775664059413
$ws.Cells.Item(20, 5).Value = 10.28517731587848
$ws.Cells.Item(20, 6).Value = 50.42786127937489
$ws.Cells.Item(20, 7).Value = 3.737047070456642
$ws.Cells.Item(20, 9).Value = 36.97979171652223
$ws.Cells.Item(20, 10).Value = 10.24354753617051
$ws.Cells.Item(20, 12).Value = 11.32870186375459
$ws.Cells.Item(20, 13).Value = 18.65972607189615

$ws.Cells.Item(21, 2).Value = 22.43523810111602
$ws.Cells.Item(21, 3).Value = 13.1362757791711
$ws.Cells.Item(21, 5).Value = 10.35418358918659
$ws.Cells.Item(21, 6).Value = 50.76136689948603
$ws.Cells.Item(21, 7).Value = 3.73217331705836
$ws.Cells.Item(21, 9).Value = 37.16499482200453
$ws.Cells.Item(21, 10).Value = 10.22643607263855
$ws.Cells.Item(21, 12).Value = 11.30824636584193
$ws.Cells.Item(21, 13).Value = 18.70873901220868

$ws.Cells.Item(22, 2).Value = 22.68789322716805
$ws.Cells.Item(22, 3).Value = 13.51585599950332
$ws.Cells.Item(22, 5).Value = 10.39865583169448
$ws.Cells.Item(22, 6).Value = 50.98713091456379
$ws.Cells.Item(22, 7).Value = 3.729101088064985
$ws.Cells.Item(22, 9).Value = 37.29140568996228
$ws.Cells.Item(22, 10).Value = 10.21574211957004
$ws.Cells.Item(22, 12).Value = 11.29612781238782
$ws.Cells.Item(22, 13).Value = 18.74425157807805

$ws.Cells.Item(23, 2).Value = 22.55287581369877
$ws.Cells.Item(23, 3).Value = 13.314355793053
$ws.Cells.Item(23, 5).Value = 10.37497784485912
$ws.Cells.Item(23, 6).Value = 50.86593490045745
$ws.Cells.Item(23, 7).Value = 3.730730492808566
$ws.Cells.Item(23, 9).Value = 37.22345099660685
$ws.Cells.Item(23, 10).Value = 10.22140488366249
$ws.Cells.Item(23, 12).Value = 11.30248155239177
$ws.Cells.Item(23, 13).Value = 18.72498168813445

$ws.Cells.Item(24, 2).Value = 22.04575844140664
$ws.Cells.Item(24, 3).Value = 12.52778540531309
$ws.Cells.Item(24, 5).Value = 10.28403891620186
$ws.Cells.Item(24, 6).Value = 50.42254458986304
$ws.Cells.Item(24, 7).Value = 3.73712864744901
$ws.Cells.Item(24, 9).Value = 36.97685628924272
$ws.Cells.Item(24, 10).Value = 10.24383547119026
$ws.Cells.Item(24, 12).Value = 11.32905735961427
$ws.Cells.Item(24, 13).Value = 18.65898420706364

$ws.Cells.Item(25, 2).Value = 21.51095886510729
$ws.Cells.Item(25, 3).Value = 11.63295414855104
$ws.Cells.Item(25, 5).Value = 10.18328112344007
$ws.Cells.Item(25, 6).Value = 49.97890563740937
$ws.Cells.Item(25, 7).Value = 3.744519774412373
$ws.Cells.Item(25, 9).Value = 36.73433144341883
$ws.Cells.Item(25, 10).Value = 10.27012957718841
$ws.Cells.Item(25, 12).Value = 11.36311577325764
$ws.Cells.Item(25, 13).Value = 18.60294193790598
